$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster (column A) changed from "Resolving-Mac" to "MuSCs" for all data rows
$ws.Range("A2").Value = "MuSCs"
$ws.Range("A3").Value = "MuSCs"
$ws.Range("A4").Value = "MuSCs"
$ws.Range("A5").Value = "MuSCs"

# Row 2 (Target cluster = ECs)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.06815433333333333
$ws.Range("H2").Value = 0.204463
$ws.Range("M2").Value = 1.123319
$ws.Range("N2").Value = 3.369957
$ws.Range("O2").Value = 0.05053686506648315
$ws.Range("P2").Value = 0.05053686506648315
$ws.Range("Q2").Value = 0.07655905756566668
$ws.Range("R2").Value = 0.6890315180910001
$ws.Range("S2").Value = 0.05053686506648315
$ws.Range("T2").Value = 0.05053686506648315

# Row 3 (Target cluster = FAPs)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.06815433333333333
$ws.Range("H3").Value = 0.204463
$ws.Range("O3").Value = 0.5042195746532222
$ws.Range("P3").Value = 0.5042195746532223
$ws.Range("Q3").Value = 0.7638498231108888
$ws.Range("R3").Value = 6.874648407998
$ws.Range("S3").Value = 0.5042195746532222
$ws.Range("T3").Value = 0.5042195746532223

# Row 4 (Target cluster = MuSCs)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.06815433333333333
$ws.Range("H4").Value = 0.204463
$ws.Range("M4").Value = 4.958620666666667
$ws.Range("N4").Value = 14.875862
$ws.Range("O4").Value = 0.2230827962023326
$ws.Range("P4").Value = 0.2230827962023326
$ws.Range("Q4").Value = 0.3379514857895556
$ws.Range("R4").Value = 3.041563372106
$ws.Range("S4").Value = 0.2230827962023326
$ws.Range("T4").Value = 0.2230827962023326

# Row 5 (Target cluster = Resolving-Mac)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.06815433333333333
$ws.Range("H5").Value = 0.204463
$ws.Range("M5").Value = 4.938126
$ws.Range("N5").Value = 14.814378
$ws.Range("O5").Value = 0.222160764077962
$ws.Range("P5").Value = 0.222160764077962
$ws.Range("Q5").Value = 0.336554685446
$ws.Range("R5").Value = 3.028992169014
$ws.Range("S5").Value = 0.222160764077962
$ws.Range("T5").Value = 0.222160764077962
